# TTO 01 "Trainingskarten gestalten" - Teamtool Training card template
#
# The slide-master title placeholder ("Titelplatzhalter 1") contains two
# text runs separated by a manual line break:
#   Run 1: "Überschrift "
#   <a:br/>
#   Run 2: "bearbeiten "
#
# The edit upper-cases both pieces of text:
#   "Überschrift " -> "ÜBERSCHRIFT"
#   "bearbeiten "  -> "BEARBEITEN "
#
# Only the existing runs' text is rewritten in place, so every other
# attribute (rPr/lang/dirty/smtClean, the <a:br>, the paragraph's
# endParaRPr, shape position, ...) stays exactly as it was.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$titleShape = $master.Shapes.Item("Titelplatzhalter 1")
$titleRange = $titleShape.TextFrame.TextRange

$titleRange.Runs(1).Text = "ÜBERSCHRIFT"
$titleRange.Runs(2).Text = "BEARBEITEN "
